$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H96").Value = 594.1429000000001
$ws.Range("I96").Value = 359.83334
$ws.Range("J96").Value = 2000
$ws.Range("K96").Value = 1079.50002
$ws.Range("L96").Value = 6000
$ws.Range("M96").Value = 293.4999800000001
$ws.Range("N96").Value = -8746

$ws.Range("H100").Value = 2606.7778
$ws.Range("I100").Value = 2701.2942
$ws.Range("K100").Value = 2701.2942
$ws.Range("M100").Value = -2160.2942

$ws.Range("H111").Value = 371.45456
$ws.Range("I111").Value = 385.4
$ws.Range("K111").Value = 1156.2
$ws.Range("M111").Value = 1910.8

$ws.Range("H129").Value = 1819.6
$ws.Range("J129").Value = 3999
$ws.Range("L129").Value = 11997
$ws.Range("N129").Value = -21997

$ws.Range("H138").Value = 4036.375
$ws.Range("I138").Value = 1855.75
$ws.Range("J138").Value = 4581.5312
$ws.Range("K138").Value = 5567.25
$ws.Range("L138").Value = 13744.5936
$ws.Range("M138").Value = -427.25
$ws.Range("N138").Value = -24024.5936

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H132").Value = 4255
$ws.Range("I132").Value = 3697
$ws.Range("K132").Value = 11091
$ws.Range("M132").Value = -8561

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 3770.353
$ws.Range("I20").Value = 1418
$ws.Range("J20").Value = 8083
$ws.Range("K20").Value = 1418
$ws.Range("L20").Value = 8083
$ws.Range("M20").Value = -1171
$ws.Range("N20").Value = -8577

$ws.Range("H80").Value = 1498.875
$ws.Range("I80").Value = 886.5
$ws.Range("K80").Value = 886.5
$ws.Range("M80").Value = 111.5

$ws.Range("H83").Value = 1498.875
$ws.Range("I83").Value = 886.5
$ws.Range("K83").Value = 4432.5
$ws.Range("M83").Value = 559.5

$ws.Range("H86").Value = 2918.9546
$ws.Range("I86").Value = 1726.0625
$ws.Range("K86").Value = 1726.0625
$ws.Range("M86").Value = -603.0625

$ws.Range("H89").Value = 2918.9546
$ws.Range("I89").Value = 1726.0625
$ws.Range("K89").Value = 8630.3125
$ws.Range("M89").Value = -3014.3125

$ws.Range("H99").Value = 34336.668
$ws.Range("I99").Value = 34336.668
$ws.Range("K99").Value = 34336.668
$ws.Range("M99").Value = -32838.668

$ws.Range("H105").Value = 2947.5
$ws.Range("I105").Value = 2947.5
$ws.Range("K105").Value = 2947.5
$ws.Range("M105").Value = -1200.5

$ws.Range("H107").Value = 2479.4
$ws.Range("I107").Value = 2487.25
$ws.Range("J107").Value = 2448
$ws.Range("K107").Value = 2487.25
$ws.Range("L107").Value = 2448
$ws.Range("M107").Value = -567.25
$ws.Range("N107").Value = -6288

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1996.1
$ws.Range("I31").Value = 1506.8334
$ws.Range("J31").Value = 2730
$ws.Range("K31").Value = 1506.8334
$ws.Range("L31").Value = 2730
$ws.Range("M31").Value = -1211.8334
$ws.Range("N31").Value = -3320

$ws.Range("H34").Value = 1996.1
$ws.Range("I34").Value = 1506.8334
$ws.Range("J34").Value = 2730
$ws.Range("K34").Value = 1506.8334
$ws.Range("L34").Value = 2730
$ws.Range("M34").Value = -1304.8334
$ws.Range("N34").Value = -3134

$ws.Range("H58").Value = 4715.6665
$ws.Range("I58").Value = 2548.5
$ws.Range("K58").Value = 2548.5
$ws.Range("M58").Value = -2345.5

$ws.Range("H64").Value = 139998
$ws.Range("J64").Value = 139998
$ws.Range("L64").Value = 139998
$ws.Range("N64").Value = -140494

$ws.Range("H67").Value = 139998
$ws.Range("J67").Value = 139998
$ws.Range("L67").Value = 139998
$ws.Range("N67").Value = -141714

$ws.Range("H93").Value = 8703.5
$ws.Range("I93").Value = 8703.5
$ws.Range("K93").Value = 8703.5
$ws.Range("M93").Value = -6831.5

$ws.Range("H132").Value = 4569.4
$ws.Range("I132").Value = 4316
$ws.Range("J132").Value = 4949.5
$ws.Range("K132").Value = 12948
$ws.Range("L132").Value = 14848.5
$ws.Range("M132").Value = -10418
$ws.Range("N132").Value = -19908.5

$ws.Range("H134").Value = 2970.6667
$ws.Range("I134").Value = 0
$ws.Range("K134").Value = 0
$ws.Range("M134").Value = ""

$ws.Range("H136").Value = 4715.6665
$ws.Range("I136").Value = 2548.5
$ws.Range("K136").Value = 7645.5
$ws.Range("M136").Value = -5095.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H92").Value = 3794
$ws.Range("J92").Value = 10003
$ws.Range("L92").Value = 30009
$ws.Range("N92").Value = -32505

$ws.Range("H141").Value = 6599.4
$ws.Range("I141").Value = 6599.4
$ws.Range("K141").Value = 19798.2
$ws.Range("M141").Value = -14618.2

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 10000
$ws.Range("I70").Value = 10000
$ws.Range("K70").Value = 10000
$ws.Range("M70").Value = -9730

$ws.Range("H73").Value = 10000
$ws.Range("I73").Value = 10000
$ws.Range("K73").Value = 10000
$ws.Range("M73").Value = -9064

$ws.Range("H136").Value = 42073.5
$ws.Range("J136").Value = 42073.5
$ws.Range("L136").Value = 126220.5
$ws.Range("N136").Value = -131320.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 5499.5

$ws.Range("H136").Value = 3774.4
$ws.Range("J136").Value = 3900
$ws.Range("L136").Value = 11700
$ws.Range("N136").Value = -16800

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 552.9231
$ws.Range("I107").Value = 498.14285
$ws.Range("J107").Value = 616.8333
$ws.Range("K107").Value = 1494.42855
$ws.Range("L107").Value = 1850.4999
$ws.Range("M107").Value = 425.5714499999999
$ws.Range("N107").Value = -5690.4999

$ws.Range("H122").Value = 1571.2858
$ws.Range("I122").Value = 1625
$ws.Range("J122").Value = 1499.6666
$ws.Range("K122").Value = 4875
$ws.Range("L122").Value = 4498.9998
$ws.Range("M122").Value = -2425
$ws.Range("N122").Value = -9398.9998

$ws.Range("H136").Value = 6209.0835
$ws.Range("I136").Value = 6189.375
$ws.Range("J136").Value = 6248.5
$ws.Range("K136").Value = 18568.125
$ws.Range("L136").Value = 18745.5
$ws.Range("M136").Value = -16018.125
$ws.Range("N136").Value = -23845.5
